$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 199
$ws.Range("I5").Value = 199
$ws.Range("K5").Value = 199
$ws.Range("M5").Value = -84
$ws.Range("H6").Value = 1786.7142
$ws.Range("I6").Value = 241.4
$ws.Range("K6").Value = 724.2
$ws.Range("M6").Value = -612.2
$ws.Range("H33").Value = 69433.75999999999
$ws.Range("I33").Value = 204.13043
$ws.Range("J33").Value = 334814
$ws.Range("K33").Value = 204.13043
$ws.Range("L33").Value = 334814
$ws.Range("M33").Value = 24.86957000000001
$ws.Range("N33").Value = -335272
$ws.Range("H34").Value = 13573.375
$ws.Range("I34").Value = 9764.666999999999
$ws.Range("K34").Value = 9764.666999999999
$ws.Range("M34").Value = -9561.666999999999
$ws.Range("H36").Value = 13573.375
$ws.Range("I36").Value = 9764.666999999999
$ws.Range("K36").Value = 9764.666999999999
$ws.Range("M36").Value = -9049.666999999999
$ws.Range("H39").Value = 4153.4
$ws.Range("I39").Value = 1882
$ws.Range("J39").Value = 5222.294
$ws.Range("K39").Value = 5646
$ws.Range("L39").Value = 15666.882
$ws.Range("M39").Value = -5350
$ws.Range("N39").Value = -16258.882
$ws.Range("H81").Value = 90000
$ws.Range("J81").Value = 90000
$ws.Range("L81").Value = 90000
$ws.Range("N81").Value = -91996
$ws.Range("H84").Value = 90000
$ws.Range("J84").Value = 90000
$ws.Range("L84").Value = 270000
$ws.Range("N84").Value = -279984
$ws.Range("H125").Value = 8010.3335
$ws.Range("I125").Value = 8010.3335
$ws.Range("K125").Value = 72093.0015
$ws.Range("M125").Value = -69633.0015
$ws.Range("H132").Value = 3977.8125
$ws.Range("I132").Value = 3561.724
$ws.Range("K132").Value = 10685.172
$ws.Range("M132").Value = -8155.172

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2810.9167
$ws.Range("I2").Value = 2216.889
$ws.Range("J2").Value = 4593
$ws.Range("K2").Value = 2216.889
$ws.Range("L2").Value = 4593
$ws.Range("M2").Value = -2103.889
$ws.Range("N2").Value = -4819
$ws.Range("H5").Value = 201.42857
$ws.Range("H43").Value = 36895
$ws.Range("I43").Value = 36895
$ws.Range("K43").Value = 36895
$ws.Range("M43").Value = -36582
$ws.Range("H45").Value = 3401.8167
$ws.Range("J45").Value = 3808.0466
$ws.Range("L45").Value = 3808.0466
$ws.Range("N45").Value = -4562.0466
$ws.Range("H61").Value = 4081
$ws.Range("I61").Value = 2947.6758
$ws.Range("K61").Value = 2947.6758
$ws.Range("M61").Value = -2735.6758
$ws.Range("H88").Value = 6747
$ws.Range("I88").Value = 5999.5
$ws.Range("K88").Value = 5999.5
$ws.Range("M88").Value = -5593.5
$ws.Range("H91").Value = 6747
$ws.Range("I91").Value = 5999.5
$ws.Range("K91").Value = 5999.5
$ws.Range("M91").Value = -4595.5
$ws.Range("H116").Value = 2810.9167
$ws.Range("I116").Value = 2216.889
$ws.Range("J116").Value = 4593
$ws.Range("K116").Value = 2216.889
$ws.Range("L116").Value = 4593
$ws.Range("M116").Value = 77.11099999999988
$ws.Range("N116").Value = -9181
$ws.Range("H136").Value = 4081
$ws.Range("I136").Value = 2947.6758
$ws.Range("K136").Value = 8843.027399999999
$ws.Range("M136").Value = -6293.027399999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2810.9167
$ws.Range("I3").Value = 2216.889
$ws.Range("J3").Value = 4593
$ws.Range("K3").Value = 2216.889
$ws.Range("L3").Value = 4593
$ws.Range("M3").Value = -2102.889
$ws.Range("N3").Value = -4821
$ws.Range("H4").Value = 201.42857
$ws.Range("H22").Value = 449.33334
$ws.Range("I22").Value = 499
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 499
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -326
$ws.Range("N22").Value = -696
$ws.Range("H107").Value = 1538.6923
$ws.Range("I107").Value = 1364.8182
$ws.Range("J107").Value = 2495
$ws.Range("K107").Value = 1364.8182
$ws.Range("L107").Value = 2495
$ws.Range("M107").Value = 555.1818000000001
$ws.Range("N107").Value = -6335

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 15625812
$ws.Range("I107").Value = 22728192
$ws.Range("K107").Value = 22728192
$ws.Range("M107").Value = -22726272
$ws.Range("H141").Value = 327857
$ws.Range("J141").Value = 327857
$ws.Range("L141").Value = 327857
$ws.Range("N141").Value = -338217

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 823.4286
$ws.Range("I2").Value = 1264.9166
$ws.Range("K2").Value = 7589.4996
$ws.Range("M2").Value = -7476.4996
$ws.Range("H12").Value = 223.41667
$ws.Range("I12").Value = 8.25
$ws.Range("K12").Value = 24.75
$ws.Range("M12").Value = 148.25
$ws.Range("H46").Value = 34876.734
$ws.Range("I46").Value = 1386.7059
$ws.Range("K46").Value = 4160.1177
$ws.Range("M46").Value = -4069.1177
$ws.Range("H75").Value = 1750
$ws.Range("J75").Value = 1750
$ws.Range("L75").Value = 5250
$ws.Range("N75").Value = -7246
$ws.Range("H78").Value = 1750
$ws.Range("J78").Value = 1750
$ws.Range("L78").Value = 15750
$ws.Range("N78").Value = -25734
$ws.Range("H92").Value = 738
$ws.Range("J92").Value = 738
$ws.Range("L92").Value = 2214
$ws.Range("N92").Value = -4710

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 784
$ws.Range("I97").Value = 784.7692
$ws.Range("J97").Value = 782.3333
$ws.Range("K97").Value = 784.7692
$ws.Range("L97").Value = 782.3333
$ws.Range("M97").Value = -288.7692
$ws.Range("N97").Value = -1774.3333
$ws.Range("H107").Value = 949.3333
$ws.Range("I107").Value = 764.4286
$ws.Range("K107").Value = 764.4286
$ws.Range("M107").Value = 1155.5714
$ws.Range("H113").Value = 2941.7778
$ws.Range("I113").Value = 2081
$ws.Range("K113").Value = 2081
$ws.Range("M113").Value = 89
$ws.Range("H122").Value = 3188.4348
$ws.Range("I122").Value = 2780.6875
$ws.Range("J122").Value = 4120.4287
$ws.Range("K122").Value = 8342.0625
$ws.Range("L122").Value = 12361.2861
$ws.Range("M122").Value = -5892.0625
$ws.Range("N122").Value = -17261.2861
$ws.Range("H126").Value = 4930.091
$ws.Range("I126").Value = 6559.875
$ws.Range("K126").Value = 19679.625
$ws.Range("M126").Value = -17209.625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2525.4092
$ws.Range("I40").Value = 2117.9412
$ws.Range("K40").Value = 2117.9412
$ws.Range("M40").Value = -1981.9412
$ws.Range("H46").Value = 1755.5714
$ws.Range("I46").Value = 1848.1666
$ws.Range("K46").Value = 1848.1666
$ws.Range("M46").Value = -1660.1666
$ws.Range("H55").Value = 3505.4443
$ws.Range("I55").Value = 3158.3333
$ws.Range("J55").Value = 4199.6665
$ws.Range("K55").Value = 3158.3333
$ws.Range("L55").Value = 4199.6665
$ws.Range("M55").Value = -2985.3333
$ws.Range("N55").Value = -4545.6665
$ws.Range("H122").Value = 8848.5
$ws.Range("I122").Value = 9105.200000000001
$ws.Range("J122").Value = 4998
$ws.Range("K122").Value = 27315.6
$ws.Range("L122").Value = 14994
$ws.Range("M122").Value = -24865.6
$ws.Range("N122").Value = -19894

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 353341.66
$ws.Range("I32").Value = 353341.66
$ws.Range("K32").Value = 353341.66
$ws.Range("M32").Value = -353024.66
$ws.Range("H62").Value = 4033.3333
$ws.Range("I62").Value = 3550
$ws.Range("K62").Value = 3550
$ws.Range("M62").Value = -2926
$ws.Range("H65").Value = 4033.3333
$ws.Range("I65").Value = 3550
$ws.Range("K65").Value = 17750
$ws.Range("M65").Value = -14630
$ws.Range("H80").Value = 65000
$ws.Range("J80").Value = 65000
$ws.Range("L80").Value = 65000
$ws.Range("N80").Value = -66996
$ws.Range("H83").Value = 65000
$ws.Range("J83").Value = 65000
$ws.Range("L83").Value = 195000
$ws.Range("N83").Value = -204984
$ws.Range("H122").Value = 3539.077
$ws.Range("I122").Value = 3703
$ws.Range("J122").Value = 3347.8333
$ws.Range("K122").Value = 11109
$ws.Range("L122").Value = 10043.4999
$ws.Range("M122").Value = -8659
$ws.Range("N122").Value = -14943.4999
$ws.Range("H126").Value = 2850.25
$ws.Range("I126").Value = 2842.9473
$ws.Range("J126").Value = 2989
$ws.Range("K126").Value = 2850.25
$ws.Range("L126").Value = 8967
$ws.Range("M126").Value = -6058.841899999999
$ws.Range("N126").Value = -13907
